$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dimension string used for rows 6 & 7 ("18x6.0-10") was originally typed
# with a Unicode EN DASH (U+2013) instead of a plain hyphen-minus. Fix it so
# it matches the rest of the sheet (and drop the now-unused shared string).
$ws.Range("C6").Value = "18x6.0-10"
$ws.Range("C7").Value = "18x6.0-10"

# Update the window/selection state: scroll back to the top-left of the
# sheet and move the active cell/selection to D10 (previously the view was
# scrolled to A6 with E28 selected).
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("D10").Select()
